$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.918.04"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.632.55"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'215.91"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "'0.5105"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.2576"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "'0.06349"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").Value = "'19.48"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").Value = "'0.07772"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "'4.277"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "1.634.99"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "1.857.47"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "'0.5493"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").Value = "'63.95"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "0.0₅7659"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "25.938.29"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "'194.86"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "'4.415"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "'9.867"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'6.053"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'1.895"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").Value = "'0.1256"
$ws.Range("E27").Value = "  +4.94%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'15.60"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'6.756"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "'0.04893"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").Value = "'3.245"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").Value = "'3.190"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "'2.368"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("D38").Value = "'2.537"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("D39").Value = "1.117.21"
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("D40").Value = "'0.01557"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "'5.604"
$ws.Range("E42").Value = "  +2.74%  "
$ws.Range("D43").Value = "'0.7958"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D44").Value = "'97.53"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("D45").Value = "1.765.98"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("E46").Value = "  -7.72%  "
$ws.Range("D47").Value = "'0.4446"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").Value = "'54.76"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").Value = "'0.05129"
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("D51").Value = "'7.594"
$ws.Range("E51").Value = "  +3.20%  "
